$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value2 = 3
$ws.Cells.Item(2, 7).Value2 = 3.4
$ws.Cells.Item(2, 8).Value2 = 2.32
$ws.Cells.Item(2, 9).Value2 = 2.62
$ws.Cells.Item(2, 10).Value2 = 3.3
$ws.Cells.Item(2, 11).Value2 = 3.85
$ws.Cells.Item(2, 12).Value2 = 1.31
$ws.Cells.Item(2, 13).Value2 = 1.07
$ws.Cells.Item(2, 14).Value2 = 3.55
$ws.Cells.Item(2, 15).Value2 = 1.31
$ws.Cells.Item(2, 16).Value2 = 1.89
$ws.Cells.Item(2, 17).Value2 = 1.91
$ws.Cells.Item(2, 18).Value2 = 1.34
$ws.Cells.Item(2, 19).Value2 = 3.3
$ws.Cells.Item(2, 20).Value2 = 1.72
$ws.Cells.Item(2, 21).Value2 = 2.12
$ws.Cells.Item(2, 22).Value2 = 1.61
$ws.Cells.Item(2, 23).Value2 = 1.41
$ws.Cells.Item(2, 24).Value2 = 16.5
$ws.Cells.Item(2, 25).Value2 = 12.5
$ws.Cells.Item(2, 26).Value2 = 20
$ws.Cells.Item(2, 28).Value2 = 15.5
$ws.Cells.Item(2, 29).Value2 = 9.4
$ws.Cells.Item(2, 30).Value2 = 14
$ws.Cells.Item(2, 33).Value2 = 17
$ws.Cells.Item(2, 34).Value2 = 21

# Row 3
$ws.Cells.Item(3, 6).Value2 = 1.85
$ws.Cells.Item(3, 7).Value2 = 1.89
$ws.Cells.Item(3, 8).Value2 = 5.1
$ws.Cells.Item(3, 10).Value2 = 3.35
$ws.Cells.Item(3, 11).Value2 = 3.55
$ws.Cells.Item(3, 14).Value2 = 2.72
$ws.Cells.Item(3, 16).Value2 = 1.58
$ws.Cells.Item(3, 17).Value2 = 2.46
$ws.Cells.Item(3, 19).Value2 = 5
$ws.Cells.Item(3, 20).Value2 = 2.14
$ws.Cells.Item(3, 21).Value2 = 1.73
$ws.Cells.Item(3, 23).Value2 = 2.12
$ws.Cells.Item(3, 29).Value2 = 8
$ws.Cells.Item(3, 30).Value2 = 980

# Row 4
$ws.Cells.Item(4, 6).Value2 = 3.8
$ws.Cells.Item(4, 7).Value2 = 5.8
$ws.Cells.Item(4, 8).Value2 = 1.84
$ws.Cells.Item(4, 9).Value2 = 2.06
$ws.Cells.Item(4, 11).Value2 = 5.1
$ws.Cells.Item(4, 12).Value2 = 1.3
$ws.Cells.Item(4, 14).Value2 = 3.3
$ws.Cells.Item(4, 16).Value2 = 1.82
$ws.Cells.Item(4, 22).Value2 = 1.94

# Row 5
$ws.Cells.Item(5, 8).Value2 = 2.5
$ws.Cells.Item(5, 9).Value2 = 2.78
$ws.Cells.Item(5, 12).Value2 = 1.43
$ws.Cells.Item(5, 14).Value2 = 3.35
$ws.Cells.Item(5, 16).Value2 = 1.81
$ws.Cells.Item(5, 21).Value2 = 2.06
$ws.Cells.Item(5, 39).Value2 = 120

# Row 6
$ws.Cells.Item(6, 6).Value2 = 1.97
$ws.Cells.Item(6, 9).Value2 = 5.6
$ws.Cells.Item(6, 15).Value2 = 1.62

# Row 8
$ws.Cells.Item(8, 7).Value2 = 2.64
$ws.Cells.Item(8, 9).Value2 = 4.2
$ws.Cells.Item(8, 13).Value2 = 1.14
$ws.Cells.Item(8, 15).Value2 = 1.62

# Row 9
$ws.Cells.Item(9, 6).Value2 = 1.35
$ws.Cells.Item(9, 7).Value2 = 1.42
$ws.Cells.Item(9, 8).Value2 = 8.800000000000001
$ws.Cells.Item(9, 9).Value2 = 10.5
$ws.Cells.Item(9, 10).Value2 = 5.5
$ws.Cells.Item(9, 11).Value2 = 6.6
$ws.Cells.Item(9, 14).Value2 = 5.3
$ws.Cells.Item(9, 15).Value2 = 1.18
$ws.Cells.Item(9, 16).Value2 = 2.54
$ws.Cells.Item(9, 17).Value2 = 1.53
$ws.Cells.Item(9, 18).Value2 = 1.6
$ws.Cells.Item(9, 19).Value2 = 2.32
$ws.Cells.Item(9, 20).Value2 = 1.83
$ws.Cells.Item(9, 21).Value2 = 2.02
$ws.Cells.Item(9, 22).Value2 = 1.1
$ws.Cells.Item(9, 23).Value2 = 3.35
$ws.Cells.Item(9, 24).Value2 = 29
$ws.Cells.Item(9, 25).Value2 = 38
$ws.Cells.Item(9, 26).Value2 = 90
$ws.Cells.Item(9, 29).Value2 = 14
$ws.Cells.Item(9, 30).Value2 = 42
$ws.Cells.Item(9, 31).Value2 = 1000
$ws.Cells.Item(9, 32).Value2 = 9.800000000000001
$ws.Cells.Item(9, 34).Value2 = 26
$ws.Cells.Item(9, 35).Value2 = 120
$ws.Cells.Item(9, 40).Value2 = 5.3
$ws.Cells.Item(9, 41).Value2 = 1000

# Row 10
$ws.Cells.Item(10, 6).Value2 = 4.6
$ws.Cells.Item(10, 7).Value2 = 5.1
$ws.Cells.Item(10, 8).Value2 = 1.67
$ws.Cells.Item(10, 9).Value2 = 1.77
$ws.Cells.Item(10, 10).Value2 = 4.4
$ws.Cells.Item(10, 11).Value2 = 4.9
$ws.Cells.Item(10, 14).Value2 = 5.9
$ws.Cells.Item(10, 16).Value2 = 2.68
$ws.Cells.Item(10, 17).Value2 = 1.51
$ws.Cells.Item(10, 18).Value2 = 1.68
$ws.Cells.Item(10, 20).Value2 = 1.6
$ws.Cells.Item(10, 21).Value2 = 2.46
$ws.Cells.Item(10, 22).Value2 = 2.42
$ws.Cells.Item(10, 23).Value2 = 1.24
$ws.Cells.Item(10, 28).Value2 = 28
$ws.Cells.Item(10, 39).Value2 = 65
$ws.Cells.Item(10, 41).Value2 = 7.8

# Row 11
$ws.Cells.Item(11, 6).Value2 = 2.42
$ws.Cells.Item(11, 7).Value2 = 2.56
$ws.Cells.Item(11, 8).Value2 = 2.9
$ws.Cells.Item(11, 9).Value2 = 3.15
$ws.Cells.Item(11, 10).Value2 = 3.5
$ws.Cells.Item(11, 12).Value2 = 1.34
$ws.Cells.Item(11, 16).Value2 = 1.93
$ws.Cells.Item(11, 22).Value2 = 1.46
$ws.Cells.Item(11, 23).Value2 = 1.64
$ws.Cells.Item(11, 25).Value2 = 12.5
$ws.Cells.Item(11, 26).Value2 = 21
$ws.Cells.Item(11, 27).Value2 = 980
$ws.Cells.Item(11, 28).Value2 = 11
$ws.Cells.Item(11, 30).Value2 = 13.5
$ws.Cells.Item(11, 32).Value2 = 16.5
$ws.Cells.Item(11, 38).Value2 = 980
$ws.Cells.Item(11, 39).Value2 = 95
$ws.Cells.Item(11, 40).Value2 = 22
$ws.Cells.Item(11, 41).Value2 = 29

# Row 12
$ws.Cells.Item(12, 8).Value2 = 1.94
$ws.Cells.Item(12, 9).Value2 = 2.08
$ws.Cells.Item(12, 10).Value2 = 3.1
$ws.Cells.Item(12, 12).Value2 = 1.51
$ws.Cells.Item(12, 15).Value2 = 1.45
$ws.Cells.Item(12, 22).Value2 = 1.92
$ws.Cells.Item(12, 36).Value2 = 150

# Row 13
$ws.Cells.Item(13, 9).Value2 = 2.42
$ws.Cells.Item(13, 11).Value2 = 3.3
$ws.Cells.Item(13, 22).Value2 = 1.7

# Row 14
$ws.Cells.Item(14, 6).Value2 = 1.91
$ws.Cells.Item(14, 7).Value2 = 2.02
$ws.Cells.Item(14, 9).Value2 = 5
$ws.Cells.Item(14, 14).Value2 = 3.7
$ws.Cells.Item(14, 23).Value2 = 1.98

# Row 15
$ws.Cells.Item(15, 6).Value2 = 1.9
$ws.Cells.Item(15, 7).Value2 = 1.93
$ws.Cells.Item(15, 8).Value2 = 4
$ws.Cells.Item(15, 9).Value2 = 4.5
$ws.Cells.Item(15, 10).Value2 = 4
$ws.Cells.Item(15, 11).Value2 = 4.4
$ws.Cells.Item(15, 18).Value2 = 1.57
$ws.Cells.Item(15, 20).Value2 = 1.62
$ws.Cells.Item(15, 21).Value2 = 2.44
$ws.Cells.Item(15, 22).Value2 = 1.29
$ws.Cells.Item(15, 23).Value2 = 2.06
$ws.Cells.Item(15, 39).Value2 = 70

# Row 16
$ws.Cells.Item(16, 6).Value2 = 2.86
$ws.Cells.Item(16, 8).Value2 = 2.64
$ws.Cells.Item(16, 10).Value2 = 3.35
$ws.Cells.Item(16, 11).Value2 = 3.5
$ws.Cells.Item(16, 18).Value2 = 1.38
$ws.Cells.Item(16, 32).Value2 = 1000
$ws.Cells.Item(16, 36).Value2 = 980

# Row 17
$ws.Cells.Item(17, 7).Value2 = 1.61
$ws.Cells.Item(17, 12).Value2 = 1.37
$ws.Cells.Item(17, 15).Value2 = 1.28
$ws.Cells.Item(17, 23).Value2 = 2.62

# Row 18
$ws.Cells.Item(18, 7).Value2 = 2.4
$ws.Cells.Item(18, 14).Value2 = 4.3

# Row 19
$ws.Cells.Item(19, 10).Value2 = 3.85
$ws.Cells.Item(19, 17).Value2 = 1.76
$ws.Cells.Item(19, 36).Value2 = 22
$ws.Cells.Item(19, 37).Value2 = 19

# Row 20
$ws.Cells.Item(20, 6).Value2 = 3.7
$ws.Cells.Item(20, 7).Value2 = 3.75
$ws.Cells.Item(20, 23).Value2 = 1.36
$ws.Cells.Item(20, 28).Value2 = 13
$ws.Cells.Item(20, 29).Value2 = 7.6

# Row 21
$ws.Cells.Item(21, 9).Value2 = 2.48
$ws.Cells.Item(21, 10).Value2 = 3.2
$ws.Cells.Item(21, 11).Value2 = 3.65
$ws.Cells.Item(21, 22).Value2 = 1.67
$ws.Cells.Item(21, 23).Value2 = 1.35

# Row 22
$ws.Cells.Item(22, 6).Value2 = 1.81
$ws.Cells.Item(22, 7).Value2 = 1.98
$ws.Cells.Item(22, 8).Value2 = 4.8
$ws.Cells.Item(22, 9).Value2 = 5.6
$ws.Cells.Item(22, 10).Value2 = 3.35
$ws.Cells.Item(22, 11).Value2 = 3.9
$ws.Cells.Item(22, 12).Value2 = 1.46
$ws.Cells.Item(22, 15).Value2 = 1.39
$ws.Cells.Item(22, 16).Value2 = 1.71
$ws.Cells.Item(22, 17).Value2 = 2.14
$ws.Cells.Item(22, 21).Value2 = 1.84
$ws.Cells.Item(22, 22).Value2 = 1.21
$ws.Cells.Item(22, 23).Value2 = 2.02
$ws.Cells.Item(22, 30).Value2 = 23
$ws.Cells.Item(22, 31).Value2 = 1000
$ws.Cells.Item(22, 32).Value2 = 11
$ws.Cells.Item(22, 35).Value2 = 110
$ws.Cells.Item(22, 36).Value2 = 26
$ws.Cells.Item(22, 37).Value2 = 23
$ws.Cells.Item(22, 40).Value2 = 16.5

# Row 23
$ws.Cells.Item(23, 6).Value2 = 4.6
$ws.Cells.Item(23, 7).Value2 = 4.9
$ws.Cells.Item(23, 8).Value2 = 1.96
$ws.Cells.Item(23, 9).Value2 = 1.98
$ws.Cells.Item(23, 10).Value2 = 3.6
$ws.Cells.Item(23, 11).Value2 = 3.7
$ws.Cells.Item(23, 17).Value2 = 2.08
$ws.Cells.Item(23, 22).Value2 = 2.02
$ws.Cells.Item(23, 23).Value2 = 1.26
$ws.Cells.Item(23, 24).Value2 = 13
$ws.Cells.Item(23, 26).Value2 = 11.5
$ws.Cells.Item(23, 28).Value2 = 16
$ws.Cells.Item(23, 29).Value2 = 8
$ws.Cells.Item(23, 30).Value2 = 10.5
$ws.Cells.Item(23, 31).Value2 = 22
$ws.Cells.Item(23, 32).Value2 = 38
$ws.Cells.Item(23, 33).Value2 = 20
$ws.Cells.Item(23, 34).Value2 = 21
$ws.Cells.Item(23, 36).Value2 = 120
$ws.Cells.Item(23, 37).Value2 = 70
$ws.Cells.Item(23, 38).Value2 = 80
$ws.Cells.Item(23, 39).Value2 = 130
$ws.Cells.Item(23, 40).Value2 = 90

# Row 24
$ws.Cells.Item(24, 8).Value2 = 5.2
$ws.Cells.Item(24, 11).Value2 = 3.3

# Row 25
$ws.Cells.Item(25, 8).Value2 = 2.9
$ws.Cells.Item(25, 15).Value2 = 1.7
$ws.Cells.Item(25, 22).Value2 = 1.48
$ws.Cells.Item(25, 28).Value2 = 8
$ws.Cells.Item(25, 38).Value2 = 100
$ws.Cells.Item(25, 40).Value2 = 90

# Row 26
$ws.Cells.Item(26, 14).Value2 = 2.48
$ws.Cells.Item(26, 20).Value2 = 2.38
